$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'Dr. Hend Mahmoud, Dr. Eman Tantawi, Dr. Majorelle Magdy, Dr. Servinaz Sayed Mohammad'
$ws.Range("G3").Value = 'Dr. Eman Tantawi, Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Asmaa Reda, Dr. Veronia Rafat'
$ws.Range("G4").Value = 'Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Veronia Rafat'
$ws.Range("G5").Value = 'Dr. Mohammad El-Tanany, Dr. Eman Tantawi, Dr. Nesma, Dr. Hanan Ragab, Dr. Hend Mahmoud, Dr. Nourhan Mahmoud, Dr. Veronia Rafat, Dr. Servinaz Sayed Mohammad'
$ws.Range("G6").Value = 'Dr. Nahla Nagiub, Dr. Gehan Adel, Dr. Eman Tantawi, Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Asmaa Reda, Dr. Hend Mahmoud, Dr. Nourhan Mahmoud, Dr. Veronia Rafat, Dr. Servinaz Sayed Mohammad'
$ws.Range("G7").Value = 'Dr. Eman Tantawi, Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat'
$ws.Range("G8").Value = 'Administrator, Dr. Manar Montaser, Dr. Eman Tantawi, Dr. Shimaa Ahmad Mekki, Dr. Asmaa Reda, Dr. Majorelle Magdy'
$ws.Range("G9").Value = 'Dr. Manar Montaser, Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Majorelle Magdy, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Gehan Adel'
$ws.Range("G10").Value = 'Dr. Sara Wael, Dr. Shimaa Ahmad Mekki, Dr. Rana Abo-Zaid, Dr. Alshimaa Atef, Dr. Heba Mahmoud Ali, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad'
$ws.Range("G11").Value = 'Dr. Hend Mahmoud, Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Asmaa Reda'
$ws.Range("G13").Value = 'D Wessam Atef, Dr. Omnia Mohammad, Dr. Shimaa Ashraf, Dr. Safa Hany, Dr. Mariam Nour El-Din'
$ws.Range("G15").Value = 'Dr. Amal Awwad, D Wessam Atef'
$ws.Range("G17").Value = 'Dr. Esraa Mostafa, Dr. Marwa Mustafa, Dr. Sarah Abdelmohsen, Dr. Arwa Al-Sayed, Dr. Nourhan Osama, Dr. Dina Adel, Dr. Madeha Saeed, Dr. Basma Hamed, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya'
$ws.Range("G24").Value = 'Dr. Maryam Ashraf, Dr. Yasmin, Dr. Aya Emad, Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Youstina Magdy, Dr. Remon, Dr. Monica, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida, Dr. Salma Hassan'
$ws.Range("G25").Value = 'Dr. Aya Emad, Dr. Marina Atef, Dr. Abdullah El-Agrody, Dr. Remon, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Youstina Magdy'
$ws.Range("G27").Value = 'Dr. Yasmin, Dr. Neveen Nashaat, Dr. Eman Mohammad Al, Dr. Remon, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida, Dr. Eman Samir Gabry, Dr. Salma Hassan'
$ws.Range("G28").Value = 'Dr. Nardine, Dr. Neveen Nashaat, Dr. Abdullah El-Agrody, Dr. Aya Hanafy, Dr. Remon, Dr. Wafaa Ebida, Dr. Eman Samir Gabry, Dr. Salma Hassan'
$ws.Range("G29").Value = 'Dr. Neveen Nashaat, Dr. Naema Gomaa, Dr. Remon, Dr. Monica, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry'
$ws.Range("G30").Value = 'Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Veronia Rafat'
$ws.Range("G31").Value = 'Dr. Eman Tantawi, Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Asmaa Reda, Dr. Veronia Rafat'
$ws.Range("G32").Value = 'Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Veronia Rafat'
$ws.Range("G33").Value = 'Dr. Mohammad El-Tanany, Dr. Eman Tantawi, Dr. Nesma, Dr. Hanan Ragab, Dr. Hend Mahmoud, Dr. Nourhan Mahmoud, Dr. Veronia Rafat, Dr. Servinaz Sayed Mohammad'
$ws.Range("G34").Value = 'Dr. Nahla Nagiub, Dr. Gehan Adel, Dr. Eman Tantawi, Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Asmaa Reda, Dr. Hend Mahmoud, Dr. Nourhan Mahmoud, Dr. Veronia Rafat, Dr. Servinaz Sayed Mohammad'
$ws.Range("G35").Value = 'Dr. Eman Tantawi, Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat'
$ws.Range("G36").Value = 'Administrator, Dr. Manar Montaser, Dr. Eman Tantawi, Dr. Shimaa Ahmad Mekki, Dr. Asmaa Reda, Dr. Majorelle Magdy'
$ws.Range("G37").Value = 'Dr. Manar Montaser, Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Majorelle Magdy, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Gehan Adel'
$ws.Range("G38").Value = 'Dr. Sara Wael, Dr. Shimaa Ahmad Mekki, Dr. Rana Abo-Zaid, Dr. Alshimaa Atef, Dr. Heba Mahmoud Ali, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad'
$ws.Range("G39").Value = 'Dr. Hend Mahmoud, Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Asmaa Reda'
$ws.Range("G41").Value = 'D Wessam Atef, Dr. Omnia Mohammad, Dr. Shimaa Ashraf, Dr. Safa Hany, Dr. Mariam Nour El-Din'
$ws.Range("G43").Value = 'Dr. Amal Awwad, D Wessam Atef'
$ws.Range("G45").Value = 'Dr. Esraa Mostafa, Dr. Marwa Mustafa, Dr. Sarah Abdelmohsen, Dr. Arwa Al-Sayed, Dr. Nourhan Osama, Dr. Dina Adel, Dr. Madeha Saeed, Dr. Basma Hamed, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya'
$ws.Range("G52").Value = 'Dr. Maryam Ashraf, Dr. Yasmin, Dr. Aya Emad, Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Youstina Magdy, Dr. Remon, Dr. Monica, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida, Dr. Salma Hassan'
$ws.Range("G53").Value = 'Dr. Aya Emad, Dr. Marina Atef, Dr. Abdullah El-Agrody, Dr. Remon, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Youstina Magdy'
$ws.Range("G55").Value = 'Dr. Yasmin, Dr. Neveen Nashaat, Dr. Eman Mohammad Al, Dr. Remon, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida, Dr. Eman Samir Gabry, Dr. Salma Hassan'
$ws.Range("G56").Value = 'Dr. Nardine, Dr. Neveen Nashaat, Dr. Abdullah El-Agrody, Dr. Aya Hanafy, Dr. Remon, Dr. Wafaa Ebida, Dr. Eman Samir Gabry, Dr. Salma Hassan'
$ws.Range("G57").Value = 'Dr. Neveen Nashaat, Dr. Naema Gomaa, Dr. Remon, Dr. Monica, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry'
